$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 0.2060085836909871
$ws.Range("C2").Value = 0.5665236051502146
$ws.Range("J2").Value = 0.004291845493562232
$ws.Range("P2").Value = 0.1373390557939914
$ws.Range("S2").Value = 0.08583690987124463

# Row 3 updates
$ws.Range("B3").Value = 0.007352941176470588
$ws.Range("C3").Value = 0.02205882352941177
$ws.Range("P3").Value = 0.8602941176470589
$ws.Range("S3").Value = 0.1102941176470588

# Row 4 updates
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.6976744186046512
$ws.Range("S4").Value = 0.2790697674418605

# Row 6 updates
$ws.Range("B6").Value = 0.05882352941176471
$ws.Range("D6").Value = 0.01764705882352941
$ws.Range("F6").Value = 0.05882352941176471
$ws.Range("J6").Value = 0.2588235294117647
$ws.Range("O6").Value = 0.01764705882352941
$ws.Range("Q6").Value = 0.1764705882352941
$ws.Range("R6").Value = 0.1117647058823529
$ws.Range("S6").Value = 0.3

# Row 7 updates
$ws.Range("B7").Value = 0.1124260355029586
$ws.Range("D7").Value = 0.01775147928994083
$ws.Range("E7").Value = 0.005917159763313609
$ws.Range("F7").Value = 0.04733727810650887
$ws.Range("J7").Value = 0.1656804733727811
$ws.Range("O7").Value = 0.01183431952662722
$ws.Range("Q7").Value = 0.1420118343195266
$ws.Range("R7").Value = 0.05917159763313609
$ws.Range("S7").Value = 0.4378698224852071

# Row 8 updates
$ws.Range("B8").Value = 0.0650887573964497
$ws.Range("D8").Value = 0.02071005917159763
$ws.Range("E8").Value = 0.002958579881656805
$ws.Range("F8").Value = 0.05325443786982249
$ws.Range("J8").Value = 0.1449704142011834
$ws.Range("O8").Value = 0.01775147928994083
$ws.Range("Q8").Value = 0.1834319526627219
$ws.Range("R8").Value = 0.09171597633136094
$ws.Range("S8").Value = 0.4201183431952663

# Row 9 updates
$ws.Range("B9").Value = 0.1764705882352941
$ws.Range("D9").Value = 0.01764705882352941
$ws.Range("F9").Value = 0.04117647058823529
$ws.Range("J9").Value = 0.1176470588235294
$ws.Range("O9").Value = 0.005882352941176471
$ws.Range("Q9").Value = 0.1235294117647059
$ws.Range("R9").Value = 0.1176470588235294
$ws.Range("S9").Value = 0.4

# Row 10 updates
$ws.Range("B10").Value = 0.0888129803586678
$ws.Range("D10").Value = 0.02391118701964133
$ws.Range("E10").Value = 0.001707941929974381
$ws.Range("F10").Value = 0.05807002561912895
$ws.Range("J10").Value = 0.1093082835183604
$ws.Range("O10").Value = 0.01280956447480786
$ws.Range("Q10").Value = 0.2203245089666951
$ws.Range("R10").Value = 0.1067463706233988
$ws.Range("S10").Value = 0.3783091374893254

# Row 11 updates
$ws.Range("G11").Value = 0.1218637992831541
$ws.Range("J11").Value = 0.07526881720430108
$ws.Range("K11").Value = 0.1971326164874552
$ws.Range("L11").Value = 0.5842293906810035
$ws.Range("S11").Value = 0.02150537634408602

# Row 12 updates
$ws.Range("G12").Value = 0.686046511627907
$ws.Range("J12").Value = 0.2267441860465116
$ws.Range("K12").Value = 0.01162790697674419
$ws.Range("L12").Value = 0.04069767441860465
$ws.Range("S12").Value = 0.03488372093023256

# Row 13 updates
$ws.Range("G13").Value = 0.5555555555555556
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.1111111111111111

# Row 15 updates
$ws.Range("F15").Value = 0.025
$ws.Range("H15").Value = 0.19
$ws.Range("I15").Value = 0.045
$ws.Range("J15").Value = 0.365
$ws.Range("K15").Value = 0.09
$ws.Range("M15").Value = 0.01
$ws.Range("N15").Value = 0.005
$ws.Range("O15").Value = 0.07
$ws.Range("S15").Value = 0.2

# Row 16 updates
$ws.Range("H16").Value = 0.1494252873563219
$ws.Range("I16").Value = 0.103448275862069
$ws.Range("J16").Value = 0.4080459770114943
$ws.Range("K16").Value = 0.1149425287356322
$ws.Range("M16").Value = 0.005747126436781609
$ws.Range("O16").Value = 0.04022988505747126
$ws.Range("S16").Value = 0.1781609195402299

# Row 17 updates
$ws.Range("F17").Value = 0.02046035805626599
$ws.Range("H17").Value = 0.1611253196930946
$ws.Range("I17").Value = 0.08439897698209718
$ws.Range("J17").Value = 0.3887468030690537
$ws.Range("K17").Value = 0.1125319693094629
$ws.Range("M17").Value = 0.03069053708439898
$ws.Range("N17").Value = 0.002557544757033248
$ws.Range("O17").Value = 0.06905370843989769
$ws.Range("S17").Value = 0.1304347826086956

# Row 18 updates
$ws.Range("F18").Value = 0.01456310679611651
$ws.Range("H18").Value = 0.1407766990291262
$ws.Range("I18").Value = 0.116504854368932
$ws.Range("J18").Value = 0.4660194174757282
$ws.Range("K18").Value = 0.07766990291262135
$ws.Range("O18").Value = 0.08737864077669903
$ws.Range("S18").Value = 0.0970873786407767

# Row 19 updates
$ws.Range("F19").Value = 0.01462522851919561
$ws.Range("H19").Value = 0.1727605118829982
$ws.Range("I19").Value = 0.07769652650822668
$ws.Range("J19").Value = 0.4076782449725777
$ws.Range("K19").Value = 0.1115173674588665
$ws.Range("M19").Value = 0.01919561243144424
$ws.Range("N19").Value = 0.0009140767824497258
$ws.Range("O19").Value = 0.07861060329067641
$ws.Range("S19").Value = 0.1170018281535649

